# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -----------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.120.0.3"
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 98.7

# Totals row
$ws.Range("C4").Value = 17

# --- Good Drivers table -------------------------------------------------
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B12").Value = 56018

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B13").Value = 34244

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 442178
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2024-11-10"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B15").Value = 77849
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2021-08-18"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2020-08-05"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2019-12-14"

# Rows 18 & 19 no longer have entries - wipe them out entirely
$ws.Range("A18:E19").Clear()

# Column A got a touch narrower (engine rounds ColumnWidth to the nearest
# pixel boundary, so 43.17 is the input that comes back out as a clean 44)
$ws.Columns.Item(1).ColumnWidth = 43.17
